$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - style matches existing header cells (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for columns I (I0) and J (IF), rows 2-14
$dataI = @(8, 9, 2, 5, 7, 8, 6, 4, 8, 4, 5, 6, 2)
$dataJ = @(9, 9, 5, 8, 9, 9, 8, 7, 8, 6, 5, 6, 3)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
